$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "George Bush"
$ws.Range("B15").NumberFormat = "d-mmm-yy"
$ws.Range("B15").Value = "2/24/1988"

$ws.Range("A16").Value = "Bill Clinton"
$ws.Range("B16").NumberFormat = "d-mmm-yy"
$ws.Range("B16").Value = "2/24/1978"

$ws.Range("B16").Select()
